$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix B22: change from text "3" to a real number 3
$ws.Range("B22").Value = 3

# Add new row 23 with annotation data
$ws.Range("A23").Value = "Ying Tang"
$ws.Range("B23").Value = "3"
$ws.Range("C23").Value = "无"
$ws.Range("D23").Value = "SMY"
$ws.Range("E23").Value = "MET"
$ws.Range("F23").Value = "55e6f9d2-bdcc-4319-8467-87a8dbd0172d"
$ws.Range("G23").Value = "Byt3oJ-0W_annotated.xlsx"
$ws.Range("H23").Value = "The authors propose a new method that approximates the discrete max-weight matching by a continuous Sinkhorn operator, which looks like an analog of softmax operator on matrices."
